$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header: add D1 = "비용xxxxxxx" with the same visual style as C1 (bold
#    header style, fill + border) so it renders consistently with A1:C1.
# ---------------------------------------------------------------------------
$c1 = $ws.Cells.Item(1, 3)
$d1 = $ws.Cells.Item(1, 4)
$d1.Value = "비용xxxxxxx"
$c1.Copy()
$d1.PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) "Type checking fix": column C (비용/cost) was wrongly formatted with a
#    thousands-separator numeric format (#,##0). Clear that number format
#    back to General while preserving the existing fill/shading, using a
#    helper cell (elsewhere on the sheet) that carries the desired
#    "filled, General format" look, then paste just the formatting onto
#    C2:C5. The helper cell is removed afterwards.
# ---------------------------------------------------------------------------
$c2 = $ws.Cells.Item(2, 3)
$fillColorIndex = $c2.Interior.ColorIndex

$helper1 = $ws.Cells.Item(50, 50)
$helper1.Interior.ColorIndex = $fillColorIndex
$helper1.Copy()
$ws.Range("C2:C5").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0
$helper1.Clear()

# ---------------------------------------------------------------------------
# 3) Add the new D2:D5 numeric values (0.1 stored as float-precision double)
#    and give them the same filled/General-format look as column C.
# ---------------------------------------------------------------------------
$ws.Cells.Item(2, 4).Value = 0.10000000149011612
$ws.Cells.Item(3, 4).Value = 0.10000000149011612
$ws.Cells.Item(4, 4).Value = 0.10000000149011612
$ws.Cells.Item(5, 4).Value = 0.10000000149011612

$helper2 = $ws.Cells.Item(51, 51)
$helper2.Interior.ColorIndex = $fillColorIndex
$helper2.Copy()
$ws.Range("D2:D5").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0
$helper2.Clear()
